$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at row 3 on Tabelle1, shifting existing rows down
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "Juenger als 23 oder geboren vor 1940"

# Copy number format from B2 to B3 (same integer number format)
$ws1.Range("B2").Copy()
$ws1.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("B3").Value = "nein"

# Add a new hidden worksheet "Tabelle2" (placed after Tabelle1) with the boolean list values
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"
$ws2.Range("A1").Value = "boolean"
$ws2.Range("A2").Value = "ja"
$ws2.Range("A3").Value = "nein"
$ws2.Range("A1:A3").Select() | Out-Null

# Add data validation (list) to B3, referencing the boolean list on Tabelle2
$ws1.Range("B3").Validation.Add(3, 1, 1, "=Tabelle2!`$A`$2:`$A`$3")

# Hide Tabelle2
$ws2.Visible = 0

# Set selection on Tabelle1 to A4
$ws1.Activate() | Out-Null
$ws1.Range("A4").Select() | Out-Null
